$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.02%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.49%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.976"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.52%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07903"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.68%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.131"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-17.71%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.852"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.20%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.794"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.30%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9287"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.12%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1757"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.10%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07954"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.26%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08819"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.49%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03123"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.92%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1004"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.40%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.07%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005935"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.74%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.49%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.84%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.50%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1289"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.23%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.166"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.82%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1790"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.74%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.14%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.82%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004513"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.72%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001249"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.39%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01737"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.93%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04827"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.27%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007370"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "6.08%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1370"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.68%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002348"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.57%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01109"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "15.65%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006070"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.02%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.22%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8205"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.70%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.22%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.22%"
